$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct mapping as discussed
$ws.Range("A2").Value = "AS_NEW"
$ws.Range("B2").Value = 5

$ws.Range("A3").Value = "JUST_RENOVATED"
$ws.Range("B3").Value = 4

$ws.Range("A4").Value = "GOOD"
$ws.Range("B4").Value = 3

$ws.Range("A5").Value = "TO_RENOVATE"
$ws.Range("B5").Value = 2

$ws.Range("A6").Value = "TO_BE_DONE_UP"
$ws.Range("B6").Value = 2

$ws.Range("A7").Value = "TO_RESTORE"
$ws.Range("B7").Value = 1
